$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Helldrivers 2 cover image filename (subscript digit -> plain underscore+digit)
$ws.Range("D4").Value = "cover_helldrivers_2.jpg"

# Fix the Nioh 3 cover image filename (subscript digit -> plain underscore+digit)
$ws.Range("D5").Value = "cover_nioh_3.jpg"

# Add the newly-uploaded Nioh 3 trailer video to the video_url column
$ws.Range("E5").Value = "Nioh_3.mp4"

# Update Nioh 3 commercial data: price drop, now flagged as new, no longer on offer
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0

# Leave the selection where the editor ended up
$ws.Range("F6").Select()
